$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Índice, Distancia, max, min, Tempo)
$data = @(
    @(0, 4273.066666666667, 4649, 3812, 0.0825863758722941),
    @(1, 4389.033333333334, 4724, 3965, 0.08157630761464436),
    @(2, 4558.133333333333, 4971, 3870, 0.08572473526000976),
    @(3, 4267.033333333334, 4628, 3367, 0.08392372926076254),
    @(4, 3778.8, 4051, 3490, 0.08726345698038737),
    @(5, 3998.833333333333, 4300, 3596, 0.08494497934977213),
    @(6, 4651.2, 5005, 4186, 0.08632264931996664),
    @(7, 4125.833333333333, 4511, 3745, 0.08677010536193848),
    @(8, 4290.333333333333, 4646, 3904, 0.08525562286376953),
    @(9, 4185.9, 4576, 3737, 0.08292495409647624)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
